$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 5678
$ws.Range("E2").Value = 86
$ws.Range("F2").Value = 86
$ws.Range("G2").Value = -137
$ws.Range("H2").Value = -153
$ws.Range("I2").Value = -154
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 4084
$ws.Range("L2").Value = 2068
$ws.Range("M2").Value = 2015
$ws.Range("N2").Value = 2012
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 1053
$ws.Range("Q2").Value = 121
$ws.Range("R2").Value = -138
$ws.Range("S2").Value = 21
$ws.Range("T2").Value = 125
$ws.Range("U2").Value = -4
$ws.Range("V2").Value = 1406
$ws.Range("W2").Value = 1.52
$ws.Range("X2").Value = -2.7
$ws.Range("Y2").Value = -8.279999999999999
$ws.Range("Z2").Value = -3.69
$ws.Range("AA2").Value = 102.63
$ws.Range("AB2").Value = 120.15
$ws.Range("AC2").Value = -772
$ws.Range("AD2").Value = -17.3
$ws.Range("AE2").Value = 9072
$ws.Range("AF2").Value = 1.47
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 22869739

# Row 3
$ws.Range("D3").Value = 6523
$ws.Range("E3").Value = 165
$ws.Range("F3").Value = 165
$ws.Range("G3").Value = 304
$ws.Range("H3").Value = 286
$ws.Range("I3").Value = 279
$ws.Range("J3").Value = 6
$ws.Range("K3").Value = 5412
$ws.Range("L3").Value = 2828
$ws.Range("M3").Value = 2584
$ws.Range("N3").Value = 2491
$ws.Range("O3").Value = 93
$ws.Range("P3").Value = 1105
$ws.Range("Q3").Value = 422
$ws.Range("R3").Value = -2
$ws.Range("S3").Value = -191
$ws.Range("T3").Value = 116
$ws.Range("U3").Value = 306
$ws.Range("V3").Value = 1934
$ws.Range("W3").Value = 2.53
$ws.Range("X3").Value = 4.38
$ws.Range("Y3").Value = 12.4
$ws.Range("Z3").Value = 6.02
$ws.Range("AA3").Value = 109.42
$ws.Range("AB3").Value = 151.22
$ws.Range("AC3").Value = 1196
$ws.Range("AD3").Value = 22.06
$ws.Range("AE3").Value = 10679
$ws.Range("AF3").Value = 2.47
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 24015928

# Row 4
$ws.Range("D4").Value = 8076
$ws.Range("E4").Value = 232
$ws.Range("F4").Value = 232
$ws.Range("G4").Value = -411
$ws.Range("H4").Value = -462
$ws.Range("I4").Value = -476
$ws.Range("J4").Value = 14
$ws.Range("K4").Value = 5115
$ws.Range("L4").Value = 2963
$ws.Range("M4").Value = 2152
$ws.Range("N4").Value = 2043
$ws.Range("O4").Value = 109
$ws.Range("P4").Value = 1105
$ws.Range("Q4").Value = 279
$ws.Range("R4").Value = -81
$ws.Range("S4").Value = -68
$ws.Range("T4").Value = 156
$ws.Range("U4").Value = 123
$ws.Range("V4").Value = 1866
$ws.Range("W4").Value = 2.88
$ws.Range("X4").Value = -5.72
$ws.Range("Y4").Value = -21.01
$ws.Range("Z4").Value = -8.779999999999999
$ws.Range("AA4").Value = 137.66
$ws.Range("AB4").Value = 108.15
$ws.Range("AC4").Value = -1983
$ws.Range("AD4").Value = -7.18
$ws.Range("AE4").Value = 8760
$ws.Range("AF4").Value = 1.62
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 24015928

# Row 5
$ws.Range("D5").Value = 9353
$ws.Range("E5").Value = 305
$ws.Range("F5").Value = 305
$ws.Range("G5").Value = 214
$ws.Range("H5").Value = 187
$ws.Range("I5").Value = 187
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 5035
$ws.Range("L5").Value = 2902
$ws.Range("M5").Value = 2134
$ws.Range("N5").Value = 2131
$ws.Range("O5").Value = 3
$ws.Range("P5").Value = 1105
$ws.Range("Q5").Value = 416
$ws.Range("R5").Value = -281
$ws.Range("S5").Value = -289
$ws.Range("T5").Value = 262
$ws.Range("U5").Value = 155
$ws.Range("V5").Value = 1649
$ws.Range("W5").Value = 3.26
$ws.Range("X5").Value = 2
$ws.Range("Y5").Value = 8.960000000000001
$ws.Range("Z5").Value = 3.68
$ws.Range("AA5").Value = 136
$ws.Range("AB5").Value = 125.66
$ws.Range("AC5").Value = 779
$ws.Range("AD5").Value = 19.75
$ws.Range("AE5").Value = 9136
$ws.Range("AF5").Value = 1.68
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 24015928

# Row 6
$ws.Range("D6").Value = 8419
$ws.Range("E6").Value = 156
$ws.Range("F6").Value = 156
$ws.Range("G6").Value = -139
$ws.Range("H6").Value = -199
$ws.Range("I6").Value = -200
$ws.Range("K6").Value = 5130
$ws.Range("L6").Value = 3112
$ws.Range("M6").Value = 2019
$ws.Range("N6").Value = 2015
$ws.Range("P6").Value = 1105
$ws.Range("Q6").Value = 224
$ws.Range("R6").Value = -466
$ws.Range("S6").Value = 243
$ws.Range("T6").Value = 413
$ws.Range("U6").Value = -189
$ws.Range("V6").Value = 1901
$ws.Range("W6").Value = 1.85
$ws.Range("X6").Value = -2.37
$ws.Range("Y6").Value = -9.640000000000001
$ws.Range("Z6").Value = -3.92
$ws.Range("AA6").Value = 154.15
$ws.Range("AB6").Value = 110.2
$ws.Range("AC6").Value = -832
$ws.Range("AD6").Value = -6.88
$ws.Range("AE6").Value = 8640
$ws.Range("AF6").Value = 0.66
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 24015928

# Row 7
$ws.Range("D7").Value = 9556
$ws.Range("E7").Value = 246
$ws.Range("G7").Value = 147
$ws.Range("H7").Value = 113
$ws.Range("I7").Value = 113
$ws.Range("K7").Value = 5410
$ws.Range("L7").Value = 3279
$ws.Range("M7").Value = 2131
$ws.Range("N7").Value = 2128
$ws.Range("P7").Value = 1105
$ws.Range("Q7").Value = 337
$ws.Range("R7").Value = -293
$ws.Range("S7").Value = -80
$ws.Range("T7").Value = 250
$ws.Range("U7").Value = 87
$ws.Range("W7").Value = 2.57
$ws.Range("X7").Value = 1.18
$ws.Range("Y7").Value = 5.46
$ws.Range("Z7").Value = 2.14
$ws.Range("AA7").Value = 153.87
$ws.Range("AC7").Value = 408
$ws.Range("AD7").Value = 21.44
$ws.Range("AE7").Value = 6773
$ws.Range("AF7").Value = 1.29
$ws.Range("AG7").Value = 0

# Row 8
$ws.Range("D8").Value = 11476
$ws.Range("E8").Value = 356
$ws.Range("G8").Value = 247
$ws.Range("H8").Value = 204
$ws.Range("I8").Value = 204
$ws.Range("K8").Value = 5734
$ws.Range("L8").Value = 3399
$ws.Range("M8").Value = 2335
$ws.Range("N8").Value = 2333
$ws.Range("P8").Value = 1105
$ws.Range("Q8").Value = 514
$ws.Range("R8").Value = -402
$ws.Range("S8").Value = 0
$ws.Range("T8").Value = 360
$ws.Range("U8").Value = 154
$ws.Range("W8").Value = 3.1
$ws.Range("X8").Value = 1.78
$ws.Range("Y8").Value = 9.15
$ws.Range("Z8").Value = 3.66
$ws.Range("AA8").Value = 145.57
$ws.Range("AC8").Value = 635
$ws.Range("AD8").Value = 13.77
$ws.Range("AE8").Value = 7426
$ws.Range("AF8").Value = 1.18
$ws.Range("AG8").Value = 0

# Row 9
$ws.Range("D9").Value = 14035
$ws.Range("E9").Value = 441
$ws.Range("G9").Value = 333
$ws.Range("H9").Value = 271
$ws.Range("I9").Value = 272
$ws.Range("K9").Value = 6193
$ws.Range("L9").Value = 3587
$ws.Range("M9").Value = 2606
$ws.Range("N9").Value = 2605
$ws.Range("P9").Value = 1105
$ws.Range("Q9").Value = 511
$ws.Range("R9").Value = -443
$ws.Range("S9").Value = 0
$ws.Range("T9").Value = 400
$ws.Range("U9").Value = 111
$ws.Range("W9").Value = 3.14
$ws.Range("X9").Value = 1.93
$ws.Range("Y9").Value = 11.02
$ws.Range("Z9").Value = 4.54
$ws.Range("AA9").Value = 137.64
$ws.Range("AC9").Value = 847
$ws.Range("AD9").Value = 10.33
$ws.Range("AE9").Value = 8292
$ws.Range("AF9").Value = 1.06
$ws.Range("AG9").Value = 0

# Remove AH/AI cells for rows 7-9 (now blank/missing in target)
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()

Write-Output "Done"